# correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two "section header" rows that carried only a label in
# column A (no data in B:G): "situação do domicílio" (row 5) and
# "grandes regiões e unidades da federação" (row 8). This edit removes
# both rows outright (not just clearing the text), so every row below
# shifts up to close the gap, and the now-unused shared strings for
# those two headers are dropped from the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 ("situação do domicílio" header row, no B:G data).
$ws.Rows("5").Delete()

# After the first deletion everything shifted up by one, so the old
# row 8 ("grandes regiões e unidades da federação" header row) is now
# row 7.
$ws.Rows("7").Delete()
